{"js": "// Map of old text -> new text for every run that changes in this edit.\nconst replacements = [\n  [\"2025-08-31 Sunday\", \"2025-09-01 Monday\"],\n  [\"33\u00f72=16, 1\", \"19\u00f79=2, 1\"],\n  [\"42\u00f77=6, 0\", \"20\u00f75=4, 0\"],\n  [\"10\u00f75=2, 0\", \"62\u00f79=6, 8\"],\n  [\"89\u00f74=22, 1\", \"97\u00f77=13, 6\"],\n  [\"62\u00f77=8, 6\", \"24\u00f78=3, 0\"],\n  [\"66\u00f76=11, 0\", \"70\u00f72=35, 0\"],\n  [\"92\u00f79=10, 2\", \"82\u00f72=41, 0\"],\n  [\"57\u00f75=11, 2\", \"16\u00f78=2, 0\"],\n  [\"50\u00f73=16, 2\", \"61\u00f79=6, 7\"],\n  [\"15\u00f78=1, 7\", \"85\u00f77=12, 1\"],\n  [\"80\u00f75=16, 0\", \"79\u00f75=15, 4\"],\n  [\"96\u00f77=13, 5\", \"16\u00f75=3, 1\"],\n  [\"82\u00f75=16, 2\", \"30\u00f73=10, 0\"],\n  [\"87\u00f77=12, 3\", \"78\u00f78=9, 6\"],\n  [\"48\u00f72=24, 0\", \"97\u00f75=19, 2\"],\n  [\"11\u00f72=5, 1\", \"80\u00f79=8, 8\"],\n  [\"59\u00f72=29, 1\", \"93\u00f76=15, 3\"],\n  [\"57\u00f77=8, 1\", \"51\u00f73=17, 0\"],\n  [\"92\u00f79=10, 2\", \"91\u00f76=15, 1\"],\n  [\"98\u00f75=19, 3\", \"42\u00f74=10, 2\"],\n  [\"50\u00f76=8, 2\", \"98\u00f74=24, 2\"],\n  [\"60\u00f72=30, 0\", \"34\u00f76=5, 4\"],\n  [\"74\u00f74=18, 2\", \"73\u00f78=9, 1\"],\n  [\"28\u00f74=7, 0\", \"63\u00f73=21, 0\"],\n  [\"56\u00f72=28, 0\", \"81\u00f78=10, 1\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Gather every run-bearing paragraph in document order (title + table cells,\n// each one holding exactly one of the text values we need to replace).\nconst targets = [];\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.length > 0) {\n    targets.push(p);\n  }\n}\n\n// Replacements must be applied in document order since some old values\n// (e.g. \"92\u00f79=10, 2\") repeat more than once with different replacements.\nif (targets.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} text-bearing paragraphs, found ${targets.length}`\n  );\n}\n\nfor (let i = 0; i < targets.length; i++) {\n  targets[i].insertText(replacements[i][1], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Title paragraph date line.\n$d.Paragraphs.Item(1).Range.Text = \"2025-09-01 Monday\"\n\n# The answer table: 5 data rows (Word 1-based row numbers 1, 5, 9, 13, 17),\n# each with 5 columns. Values are set by row/column position (not by\n# Find/Replace) because one old value (\"92\u00f79=10, 2\") repeats twice with two\n# different replacements, so a global replace-all would be wrong.\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n  @(\"19\u00f79=2, 1\", \"20\u00f75=4, 0\", \"62\u00f79=6, 8\", \"97\u00f77=13, 6\", \"24\u00f78=3, 0\"),\n  @(\"70\u00f72=35, 0\", \"82\u00f72=41, 0\", \"16\u00f78=2, 0\", \"61\u00f79=6, 7\", \"85\u00f77=12, 1\"),\n  @(\"79\u00f75=15, 4\", \"16\u00f75=3, 1\", \"30\u00f73=10, 0\", \"78\u00f78=9, 6\", \"97\u00f75=19, 2\"),\n  @(\"80\u00f79=8, 8\", \"93\u00f76=15, 3\", \"51\u00f73=17, 0\", \"91\u00f76=15, 1\", \"42\u00f74=10, 2\"),\n  @(\"98\u00f74=24, 2\", \"34\u00f76=5, 4\", \"73\u00f78=9, 1\", \"63\u00f73=21, 0\", \"81\u00f78=10, 1\")\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n  $rowIndex = $dataRows[$i]\n  for ($col = 1; $col -le 5; $col++) {\n    $t.Cell($rowIndex, $col).Range.Text = $newValues[$i][$col - 1]\n  }\n}\n"}
